$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the existing _GoBack bookmark (it will be recreated at the very
#    end of the document once the new paragraphs have been appended).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Append the two new trailing paragraphs: an empty paragraph, followed by
#    a paragraph containing "mas..". A temporary placeholder run is appended
#    right after so that subsequent Range/Bookmark operations never need to
#    touch the absolute last character position of the document.
# ---------------------------------------------------------------------------
$rEnd = $d.Content
$rEnd.Collapse(0)
$rEnd.InsertParagraphAfter()

$rEnd2 = $d.Content
$rEnd2.Collapse(0)
$rEnd2.InsertParagraphAfter()

$rEnd3 = $d.Content
$rEnd3.Collapse(0)
$rEnd3.InsertAfter("mas..")

$rEnd4 = $d.Content
$rEnd4.Collapse(0)
$rEnd4.InsertAfter("ZZPLACEHOLDERZZ")

# ---------------------------------------------------------------------------
# 3. Re-create the _GoBack bookmark, collapsed immediately after "mas.." and
#    before the placeholder text.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$lastPara = $paras.Item($paras.Count)
$bmPos = $lastPara.Range.End - "ZZPLACEHOLDERZZ".Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 4. Remove the temporary placeholder text using Find/Replace (avoids direct
#    Range.Delete at the very end of the document).
# ---------------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("ZZPLACEHOLDERZZ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. Touch-up the w:proofErr markers on the two paragraphs whose proofing
#    state changed as a side effect of the edit (done last, after the
#    document has its final paragraph count, so neither target paragraph is
#    the very last paragraph in the body).
# ---------------------------------------------------------------------------

# "git init" paragraph: drop the spellStart/spellEnd pair that used to wrap
# "git" (only the grammar-check pair remains around it).
$pGitInit = $d.Paragraphs.Item(2)
$xmlGitInit = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p><w:pPr><w:ind w:left="707"/></w:pPr>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>git</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>init</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pGitInit.Range.InsertXML($xmlGitInit)

# "git log" paragraph: add a gramStart/gramEnd pair around "git" in addition
# to the existing spellStart/spellEnd pair.
$pGitLog = $d.Paragraphs.Item(9)
$xmlGitLog = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>git</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> log</w:t></w:r>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pGitLog.Range.InsertXML($xmlGitLog)
